# Fixing matricula of Matc65: update column A (matricula) values for rows 14-39.
# Each cell keeps its original text representation (no leading apostrophe / number
# conversion), matching the inlineStr text type used in the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newMatriculas = @{
    14 = "217216526"
    15 = "216117974"
    16 = "221117463"
    17 = "217125254"
    18 = "219218129"
    19 = "218215397"
    20 = "220117282"
    21 = "219217429"
    22 = "216216087"
    23 = "220121412"
    24 = "210201260"
    25 = "201520233"
    26 = "217117994"
    27 = "219118481"
    28 = "221119218"
    29 = "219215012"
    30 = "219121541"
    31 = "214007731"
    32 = "219215013"
    33 = "220117290"
    34 = "219118473"
    35 = "220117273"
    36 = "220120071"
    37 = "221216783"
    38 = "214120645"
    39 = "220217140"
}

# Ensure column A keeps being treated as text, so the numeric-looking
# matricula strings are not turned into actual numbers.
$ws.Range("A14:A39").NumberFormat = "@"

foreach ($row in $newMatriculas.Keys) {
    $ws.Range("A$row").Value = $newMatriculas[$row]
}
